# Updated odds figures for the 2024-11-17 FlashScore weekly fixtures sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("I4").Value = 3.1
$ws.Range("N4").Value = 7.5
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57

# Row 6
$ws.Range("G6").Value = 1.95
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 3.65
$ws.Range("J6").Value = 2.52
$ws.Range("L6").Value = 4.1
$ws.Range("O6").Value = 1.35
$ws.Range("P6").Value = 2.72
$ws.Range("Q6").Value = 2.02
$ws.Range("R6").Value = 1.62
$ws.Range("T6").Value = 2.52
$ws.Range("U6").Value = 1.87
$ws.Range("W6").Value = 6.4
$ws.Range("X6").Value = 8.5
$ws.Range("Y6").Value = 8.75
$ws.Range("Z6").Value = 16.5
$ws.Range("AA6").Value = 17
$ws.Range("AB6").Value = 32
$ws.Range("AD6").Value = 6.5
$ws.Range("AG6").Value = 9.5
$ws.Range("AH6").Value = 19
$ws.Range("AI6").Value = 12.5
$ws.Range("AJ6").Value = 55
$ws.Range("AK6").Value = 35
$ws.Range("AL6").Value = 45
$ws.Range("AN6").Value = 3.7
$ws.Range("AO6").Value = 9.75
$ws.Range("AP6").Value = 19.5
$ws.Range("AQ6").Value = 37
$ws.Range("AR6").Value = 75
$ws.Range("AU6").Value = 7.4
$ws.Range("AW6").Value = 5.4
$ws.Range("AX6").Value = 20
$ws.Range("AY6").Value = 28
$ws.Range("AZ6").Value = 110
$ws.Range("BA6").Value = 150
$ws.Range("BB6").Value = 350

# Row 7
$ws.Range("I7").Value = 3.3
$ws.Range("O7").Value = 1.18
$ws.Range("P7").Value = 4.5
$ws.Range("X7").Value = 12
$ws.Range("AD7").Value = 7

# Row 10
$ws.Range("G10").Value = 3.4
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 4.33
$ws.Range("N10").Value = 9
$ws.Range("AG10").Value = 6.5
$ws.Range("AH10").Value = 9
$ws.Range("AO10").Value = 21
$ws.Range("AW10").Value = 4
$ws.Range("AX10").Value = 12

# Row 12
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73

# Row 14
$ws.Range("G14").Value = 1.8
$ws.Range("H14").Value = 3.7
$ws.Range("I14").Value = 4.1
$ws.Range("J14").Value = 2.38
$ws.Range("U14").Value = 1.53
$ws.Range("V14").Value = 2.38
$ws.Range("AF14").Value = 34
$ws.Range("AH14").Value = 23
$ws.Range("AL14").Value = 29
$ws.Range("AO14").Value = 9.5
$ws.Range("AQ14").Value = 29
$ws.Range("BB14").Value = 126

# Row 15
$ws.Range("I15").Value = 2.6
$ws.Range("J15").Value = 4
$ws.Range("AG15").Value = 6
$ws.Range("AW15").Value = 4.33

# Row 17
$ws.Range("G17").Value = 2.88
$ws.Range("I17").Value = 2.5
$ws.Range("AH17").Value = 12
$ws.Range("AV17").Value = 51
$ws.Range("AW17").Value = 4.5
